# Applies the FFXIV leve-profit market-data refresh captured in the commit diff.
# For each changed row, write the new currentAveragePrice* / LevePrice* / LeveProfit* values
# (columns H-N) on the appropriate class sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 286.25
$ws.Range("I4").Value = 148
$ws.Range("J4").Value = 701
$ws.Range("K4").Value = 148
$ws.Range("L4").Value = 701
$ws.Range("M4").Value = -34
$ws.Range("N4").Value = -929

$ws.Range("H5").Value = 112.05
$ws.Range("I5").Value = 117.21429
$ws.Range("K5").Value = 117.21429
$ws.Range("M5").Value = -2.214290000000005

$ws.Range("H28").Value = 805.63635
$ws.Range("I28").Value = 684.46155
$ws.Range("J28").Value = 1255.7142
$ws.Range("K28").Value = 684.46155
$ws.Range("L28").Value = 1255.7142
$ws.Range("M28").Value = -199.46155
$ws.Range("N28").Value = -2225.7142

$ws.Range("H106").Value = 2807.3215
$ws.Range("I106").Value = 1955.8334
$ws.Range("J106").Value = 4340
$ws.Range("K106").Value = 1955.8334
$ws.Range("L106").Value = 4340
$ws.Range("M106").Value = -1324.8334
$ws.Range("N106").Value = -5602

$ws.Range("H107").Value = 813
$ws.Range("I107").Value = 267.6316
$ws.Range("J107").Value = 1676.5
$ws.Range("K107").Value = 267.6316
$ws.Range("L107").Value = 1676.5
$ws.Range("M107").Value = 1652.3684
$ws.Range("N107").Value = -5516.5

$ws.Range("H129").Value = 2068.3125
$ws.Range("I129").Value = 2246.1428
$ws.Range("J129").Value = 1930
$ws.Range("K129").Value = 6738.428400000001
$ws.Range("L129").Value = 5790
$ws.Range("M129").Value = -1738.428400000001
$ws.Range("N129").Value = -15790

$ws.Range("H132").Value = 2627645.5
$ws.Range("I132").Value = 3044838
$ws.Range("J132").Value = 5291.857
$ws.Range("K132").Value = 9134514
$ws.Range("L132").Value = 15875.571
$ws.Range("M132").Value = -9131984
$ws.Range("N132").Value = -20935.571

$ws.Range("H138").Value = 4841.955
$ws.Range("I138").Value = 1686.6842
$ws.Range("J138").Value = 6090.9165
$ws.Range("K138").Value = 5060.0526
$ws.Range("L138").Value = 18272.7495
$ws.Range("M138").Value = 79.94740000000002
$ws.Range("N138").Value = -28552.7495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26139.193
$ws.Range("I32").Value = 15534.8
$ws.Range("J32").Value = 70324.164
$ws.Range("K32").Value = 15534.8
$ws.Range("L32").Value = 70324.164
$ws.Range("M32").Value = -15247.8
$ws.Range("N32").Value = -70898.164

$ws.Range("H37").Value = 7322.5
$ws.Range("J37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1804.3704
$ws.Range("I86").Value = 1572.6857
$ws.Range("J86").Value = 2231.158
$ws.Range("K86").Value = 1572.6857
$ws.Range("L86").Value = 2231.158
$ws.Range("M86").Value = -449.6857
$ws.Range("N86").Value = -4477.157999999999

$ws.Range("H89").Value = 1804.3704
$ws.Range("I89").Value = 1572.6857
$ws.Range("J89").Value = 2231.158
$ws.Range("K89").Value = 7863.4285
$ws.Range("L89").Value = 11155.79
$ws.Range("M89").Value = -2247.4285
$ws.Range("N89").Value = -22387.79

$ws.Range("H134").Value = 35717470
$ws.Range("I134").Value = 62502224
$ws.Range("J134").Value = 4466.6665
$ws.Range("K134").Value = 187506672
$ws.Range("L134").Value = 13399.9995
$ws.Range("M134").Value = -187504137
$ws.Range("N134").Value = -18469.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 902.7778
$ws.Range("I35").Value = 902.7778
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 902.7778
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -608.7778

$ws.Range("H62").Value = 2803.6843
$ws.Range("I62").Value = 2661.25
$ws.Range("J62").Value = 2907.2727
$ws.Range("K62").Value = 2661.25
$ws.Range("L62").Value = 2907.2727
$ws.Range("M62").Value = -2037.25
$ws.Range("N62").Value = -4155.2727

$ws.Range("H65").Value = 2803.6843
$ws.Range("I65").Value = 2661.25
$ws.Range("J65").Value = 2907.2727
$ws.Range("K65").Value = 13306.25
$ws.Range("L65").Value = 14536.3635
$ws.Range("M65").Value = -10186.25
$ws.Range("N65").Value = -20776.3635

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4334.615
$ws.Range("I3").Value = 3304.5454
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 9913.636200000001
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = -9801.636200000001
$ws.Range("N3").Value = -30224

$ws.Range("H113").Value = 490.91666
$ws.Range("I113").Value = 479.2
$ws.Range("J113").Value = 499.2857
$ws.Range("K113").Value = 1437.6
$ws.Range("L113").Value = 1497.8571
$ws.Range("M113").Value = 732.4000000000001
$ws.Range("N113").Value = -5837.8571

$ws.Range("H122").Value = 832.76666
$ws.Range("J122").Value = 1497.5714
$ws.Range("L122").Value = 13478.1426
$ws.Range("N122").Value = -18378.1426

$ws.Range("H131").Value = 72674.5
$ws.Range("I131").Value = 399.5
$ws.Range("J131").Value = 78234.12
$ws.Range("K131").Value = 1198.5
$ws.Range("L131").Value = 234702.36
$ws.Range("M131").Value = 3841.5
$ws.Range("N131").Value = -244782.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 44883.332
$ws.Range("I80").Value = 2312.5
$ws.Range("J80").Value = 66168.75
$ws.Range("K80").Value = 2312.5
$ws.Range("L80").Value = 66168.75
$ws.Range("M80").Value = -1314.5
$ws.Range("N80").Value = -68164.75

$ws.Range("H83").Value = 44883.332
$ws.Range("I83").Value = 2312.5
$ws.Range("J83").Value = 66168.75
$ws.Range("K83").Value = 11562.5
$ws.Range("L83").Value = 330843.75
$ws.Range("M83").Value = -6570.5
$ws.Range("N83").Value = -340827.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 4031
$ws.Range("I35").Value = 4031
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4031
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -3695

$ws.Range("H61").Value = 2495.6538
$ws.Range("I61").Value = 1249.4286
$ws.Range("J61").Value = 3949.5833
$ws.Range("K61").Value = 1249.4286
$ws.Range("L61").Value = 3949.5833
$ws.Range("M61").Value = -1047.4286
$ws.Range("N61").Value = -4353.5833

$ws.Range("H68").Value = 9992.615
$ws.Range("I68").Value = 15343.429
$ws.Range("J68").Value = 3750
$ws.Range("K68").Value = 15343.429
$ws.Range("L68").Value = 3750
$ws.Range("M68").Value = -14594.429
$ws.Range("N68").Value = -5248

$ws.Range("H71").Value = 9992.615
$ws.Range("I71").Value = 15343.429
$ws.Range("J71").Value = 3750
$ws.Range("K71").Value = 76717.145
$ws.Range("L71").Value = 18750
$ws.Range("M71").Value = -72973.145
$ws.Range("N71").Value = -26238

$ws.Range("H82").Value = 2713.125
$ws.Range("I82").Value = 1800.6666
$ws.Range("J82").Value = 3260.6
$ws.Range("K82").Value = 1800.6666
$ws.Range("L82").Value = 3260.6
$ws.Range("M82").Value = -1439.6666
$ws.Range("N82").Value = -3982.6

$ws.Range("H85").Value = 2713.125
$ws.Range("I85").Value = 1800.6666
$ws.Range("J85").Value = 3260.6
$ws.Range("K85").Value = 1800.6666
$ws.Range("L85").Value = 3260.6
$ws.Range("M85").Value = -552.6666
$ws.Range("N85").Value = -5756.6

$ws.Range("H93").Value = 5668.8076
$ws.Range("I93").Value = 6452.0527
$ws.Range("J93").Value = 3542.8572
$ws.Range("K93").Value = 6452.0527
$ws.Range("L93").Value = 3542.8572
$ws.Range("M93").Value = -5204.0527
$ws.Range("N93").Value = -6038.8572

$ws.Range("H113").Value = 2495.6538
$ws.Range("I113").Value = 1249.4286
$ws.Range("J113").Value = 3949.5833
$ws.Range("K113").Value = 1249.4286
$ws.Range("L113").Value = 3949.5833
$ws.Range("M113").Value = 920.5714
$ws.Range("N113").Value = -8289.5833

$ws.Range("H132").Value = 1867.1522
$ws.Range("I132").Value = 1280.2778
$ws.Range("K132").Value = 3840.8334
$ws.Range("M132").Value = -1310.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 831.6087
$ws.Range("I122").Value = 625.0952
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 1875.2856
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = 574.7144000000001
$ws.Range("N122").Value = -13900

$ws.Range("H136").Value = 3926.0852
$ws.Range("I136").Value = 1262.5758
$ws.Range("J136").Value = 10204.357
$ws.Range("K136").Value = 3787.7274
$ws.Range("L136").Value = 30613.071
$ws.Range("M136").Value = -1237.7274
$ws.Range("N136").Value = -35713.071
